$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item(1)

# Version: 1.0.1 -> 0.0.0
$ws1.Range("B3").Value = "0.0.0"

# Title: "ValueSet of Gender Identity" -> "Gender Identity"
$ws1.Range("B5").Value = "Gender Identity"

# Experimental: (was blank) -> false
# A leading apostrophe forces Excel to store this as literal text
# instead of auto-converting it to a Boolean value.
$ws1.Range("B7").Value = "'false"

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-11T13:00:00-03:00
$ws1.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# --- Sheet 2: "Include from LOINC" ---
$ws2 = $wb.Worksheets.Item(2)

# Remove the "LA76696-4: / Non-binary" concept row entirely (row 4),
# shifting subsequent rows up.
$ws2.Rows.Item(4).Delete()
